# Auto-generated Excel COM-interop script
# Applies translation strings for eDWIN workbook (release 2024-03)

$wb = $excel.ActiveWorkbook
$nl = [char]10

# --- New translated strings (text content only; final shared-string index
#     order is determined by the order cells are written below) ---
$s24 = 'eDWIN \"Virtual Farm\" giver brugere i Polen mulighed for at indhente, indsamle og dele oplysninger om forekomsten af ​​skadedyr i et givet område og giver meddelelser om mulige trusler i marken.' + $nl + 'eDWIN-platformen giver også adgang til data fra omkring 600 meteorologiske stationer i hele Polen, overvågning (blandt andet) temperatur, luftfugtighed, nedbørsmængde og intensitet, atmosfærisk tryk og vindhastighed og -retning.' + $nl + 'eDWIN-rådgivningsplatformen blev oprettet som en del af projektet \"Internetplatform for rådgivning og beslutningsstøtte inden for integreret plantebeskyttelse\".' + $nl + 'Platformen er helt gratis og tilgængelig for alle på computere og som en applikation på mobile enheder, men kun i øjeblikket tilgængelig for brugere i Polen.' + $nl + 'https://www.edwin.gov.pl/euslugi/wirtualne-gospodarstwo'
$s25 = 'eDWIN Platformen (kun i Polen)'
$s26 = 'Beslutningsstøtteplatform for skadedyr tilgængelig i Polen'
$s27 = 'Der eDWIN ("virtuelle Bauernhof") ermöglicht es den Nutzern in Polen, Informationen über das Auftreten von Schaderregern in einem bestimmten Gebiet zu erhalten, zu sammeln und auszutauschen, und liefert Benachrichtigungen über mögliche Bedrohungen vor Ort.' + $nl + 'Die eDWIN-Plattform bietet auch Zugang zu den Daten von etwa 600 Wetterstationen in ganz Polen, die (unter anderem) Temperatur, Luftfeuchtigkeit, Niederschlagssumme und -intensität, Luftdruck sowie Windgeschwindigkeit und -richtung überwachen.' + $nl + 'Die eDWIN-Beratungsplattform wurde im Rahmen des Projekts \"Internet Platform for Advisory and Decision Support in Integrated Plant Protection\" geschaffen. ' + $nl + 'Die Plattform ist völlig kostenlos und für jedermann auf Computern und als Anwendung auf mobilen Geräten verfügbar, aber derzeit nur für Nutzer in Polen zugänglich. ' + $nl + 'https://www.edwin.gov.pl/euslugi/wirtualne-gospodarstwo'
$s28 = 'eDWIN-Plattform (nur in Polen)'
$s29 = 'Plattform zur Unterstützung von Entscheidungen gegen Schaderreger in Polen verfügbar'
$s30 = 'Met het eDWIN-platform kunnen gebruikers in Polen informatie verkrijgen, verzamelen en delen over het voorkomen van plagen in een bepaald gebied en meldingen doen over mogelijke bedreigingen in het veld.' + $nl + 'Het eDWIN-platform biedt ook toegang tot gegevens van ongeveer 600 meteorologische stations in heel Polen, die (onder meer) temperatuur, luchtvochtigheid, neerslaghoeveelheid en -intensiteit, luchtdruk en windsnelheid en -richting monitoren.' + $nl + 'Het adviesplatform eDWIN is opgericht als onderdeel van het project "Internet Platform for Advisory and Decision Support in Integrated Plant Protection". ' + $nl + 'Het platform is volledig gratis en beschikbaar voor iedereen op computers en als applicatie op mobiele apparaten, maar momenteel alleen toegankelijk voor gebruikers in Polen. ' + $nl + 'https://www.edwin.gov.pl/euslugi/wirtualne-gospodarstwo'
$s31 = 'eDWIN Platform (Alleen in Polen)'
$s32 = 'Beslissingsondersteunend platform voor plaagbestrijding beschikbaar in Polen'
$s33 = 'eDWIN \"Virtuell gård\" lar brukere i Polen samle inn og dele informasjon om forekomsten av skadedyr i et gitt område og gir varsler om mulige trusler i åkeren.' + $nl + 'eDWIN-plattformen gir også tilgang til data fra rundt 600 meteorologiske stasjoner i hele Polen som blant annet overvåker temperatur, luftfuktighet, nedbørsmengde og intensitet, atmosfærisk trykk og vindhastighet og vindretning.' + $nl + 'eDWIN-rådgivningsplattform ble opprettet som en del av prosjektet \"Internet Platform for Advisory and Decision Support in Integrated Plant Protection\".' + $nl + 'Plattformen er helt gratis og tilgjengelig for alle med datamaskin og som en applikasjon på mobile enheter. For øyeblikket er den kun tilgjengelig for brukere i Polen.' + $nl + 'https://www.edwin.gov.pl/euslugi/wirtualne-gospodarstwo' + $nl + ''
$s34 = 'eDWIN-plattform (bare i Polen)'
$s35 = 'Plattform i Polen som gir veiledning i beslutninger om skadegjørere.'
$s36 = 'eDWIN \"Virtual farm\" gör det möjligt för användare i Polen att få, samla in och dela information om förekomsten av skadedjur i ett visst område och ger meddelanden om möjliga hot på fältet.  eDWIN-plattformen ger också tillgång till data från cirka 600 meteorologiska stationer i hela Polen, som övervakar (bland annat) temperatur, luftfuktighet, nederbörd (totalt och intensitet), atmosfärstryck och vindhastighet och riktning.  Den rådgivande eDWIN-plattformen skapades som en del av projektet "Internetplattform för rådgivning och beslutsstöd i integrerat växtskydd".  Plattformen är helt gratis och tillgänglig för alla på datorer och som en applikation på mobila enheter, men endast för närvarande tillgänglig för användare i Polen.  https://www.edwin.gov.pl/euslugi/wirtualne-gospodarstwo'
$s37 = 'eDWIN-plattform (bara i Polen)'
$s38 = 'Plattform i Polen som ger vägledning i skadegörarbeslut'
$s39 = 'eDWIN \ "Virtual Farm \" mahdollistaa käyttäjille Puolassa hankkia, kerätä ja jakaa tietoja tuholaisten esiintymisestä tietyllä alueella ja ilmoittaa mahdollisista uhista.' + $nl + 'eDWIN-alusta tarjoaa pääsyn noin 600 meteorologisen aseman tietoihin Puolassa (muun muassa lämpötila, ilman kosteus, sademäärän kokonaismäärä ja voimakkuus, ilmakehän paine ja tuulen nopeus sekä suunta).' + $nl + 'eDWIN neuvonnallinen alusta luotiin osana projektia "Internet-alusta neuvontaan ja päätöksenteon tukemiseen integroidussa kasvinsuojelussa".' + $nl + 'Alusta on täysin ilmainen ja kaikkien saatavilla sekä tietokoneella että mobiililaitteilla, mutta toistaiseksi saatavilla vain Puolassa.' + $nl + 'https://www.edwin.gov.pl/euslugi/wirtualne-gospodarstwo'
$s40 = 'Alusta kasvintuhoojien hallinnan päätöksenteon tueksi Puolassa'
$s41 = 'eDWIN“ \ „Virtualus ūkis“ \ leidžia vartotojams Lenkijoje gauti, rinkti ir dalytis informacija apie kenkėjų paplitimą tam tikroje vietovėje ir teikia pranešimus apie galimas grėsmes lauke.' + $nl + '„eDWIN“ platforma taip pat suteikia prieigą prie maždaug 600 meteorologijos stočių  duomenų esančių visoje Lenkijoje, kuriuose, stebima temperatūra, oro drėgmė, bendras kritulių kiekis ir intensyvumas, atmosferos slėgis, vėjo greitis ir kryptis ir kt.' + $nl + '„eDWIN“ konsultacinė platforma buvo sukurta įgyvendinant projektą „Internetinė platforma konsultavimui ir sprendimų priėmimui integruotai augalų apsaugai“. ' + $nl + 'Platforma yra visiems prieinama ir visiškai nemokama naudojant kompiuteriuose ir kaip mobilią programėlę telefone, tačiau šiuo metu ja gali naudotis tik vartotojai Lenkijoje. ' + $nl + 'https://www.edwin.gov.pl/euslugi/wirtualne-gospodarstwo"'
$s42 = 'eDWIN platforma (tik Lenkijoje)'
$s43 = 'Sprendimų priėmimo pagalbos sistema skirta kenkėjams, prieinama Lenkijoje'
$s44 = 'La "ferme virtuelle" eDWIN permet aux utilisateurs polonais d''obtenir, de collecter et de partager des informations sur la présence de parasites dans une zone donnée et de recevoir des notifications sur les menaces éventuelles sur le terrain.' + $nl + 'La plateforme eDWIN permet également d''accéder aux données de quelque 600 stations météorologiques réparties sur l''ensemble du territoire polonais, qui surveillent (entre autres) la température, l''humidité de l''air, le total et l''intensité des précipitations, la pression atmosphérique ainsi que la vitesse et la direction du vent.' + $nl + 'La plateforme de conseil eDWIN a été créée dans le cadre du projet "Plateforme Internet de conseil et d''aide à la décision pour la protection intégrée des plantes". ' + $nl + 'La plateforme est entièrement gratuite et disponible pour tous sur les ordinateurs et sous forme d''application sur les appareils mobiles, mais n''est actuellement accessible qu''aux utilisateurs en Pologne. ' + $nl + 'https://www.edwin.gov.pl/euslugi/wirtualne-gospodarstwo' + $nl + ''
$s45 = 'Plate-forme eDWIN (Pologne uniquement)'
$s46 = 'Une plateforme d''aide à la décision antiparasitaire disponible en Pologne'
$s47 = '"La "fattoria virtuale" eDWIN consente agli utenti in Polonia di ottenere, raccogliere e condividere informazioni sulla presenza di parassiti in una determinata area e fornisce notifiche su possibili minacce sul campo.' + $nl + 'La piattaforma eDWIN fornisce anche l''accesso ai dati di circa 600 stazioni meteorologiche in tutta la Polonia, monitorando (tra l''altro) la temperatura, l''umidità dell''aria, il totale e l''intensità delle precipitazioni, la pressione atmosferica e la velocità e direzione del vento.' + $nl + 'La piattaforma di consulenza eDWIN è stata creata nell''ambito del progetto "Internet Platform for Advisory and Decision Support in Integrated Plant Protection". ' + $nl + 'La piattaforma è completamente gratuita e disponibile per tutti su computer e come applicazione su dispositivi mobili, ma attualmente è accessibile solo agli utenti della Polonia. ' + $nl + 'https://www.edwin.gov.pl/euslugi/wirtualne-gospodarstwo"'
$s48 = 'Piattaforma eDWIN (solo per la Polonia)'
$s49 = 'Piattaforma di supporto alle decisioni sui parassiti disponibile in Polonia'
$s50 = 'Το eDWIN \"Virtual farm\" επιτρέπει στους χρήστες στην Πολωνία να λαμβάνουν, να συλλέγουν και να μοιράζονται πληροφορίες σχετικά με την εμφάνιση παρασίτων σε μια συγκεκριμένη περιοχή και παρέχει ειδοποιήσεις σχετικά με πιθανές απειλές στο πεδίο.' + $nl + 'Η πλατφόρμα eDWIN, παρέχει επίσης πρόσβαση σε δεδομένα από περίπου 600 μετεωρολογικούς σταθμούς σε όλη την Πολωνία, που παρακολουθούν (μεταξύ άλλων) τη θερμοκρασία, την υγρασία του αέρα, το σύνολο και την ένταση της βροχόπτωσης, την ατμοσφαιρική πίεση και την ταχύτητα και κατεύθυνση του ανέμου.' + $nl + 'Η συμβουλευτική πλατφόρμα eDWIN δημιουργήθηκε στο πλαίσιο του έργου \"Internet Platform for Advisory and Decision Support in Integrated Plant Protection\.Η πλατφόρμα είναι εντελώς δωρεάν και διαθέσιμη σε όλους σε υπολογιστές και ως εφαρμογή σε κινητές συσκευές, αλλά προς το παρόν είναι προσβάσιμη μόνο σε χρήστες στην Πολωνία.  https://www.edwin.gov.pl/euslugi/wirtualne-gospodarstwo"https://www.edwin.gov.pl/euslugi/wirtualne-gospodarstwo"'
$s51 = 'Πλατφόρμα eDWIN (μόνο στην Πολωνία)'
$s52 = 'Πλατφόρμα υποστήριξης αποφάσεων για παράσιτα διαθέσιμη στην Πολωνία'

# --- Sheet 'main': row 2, fill columns C:P with the 'eDWIN' value (same as B2) ---
$wsMain = $wb.Worksheets.Item("main")
$edwinName = $wsMain.Range("B2").Value2
foreach ($col in @("C","D","E","F","G","H","I","J","K","L","M","N","O","P")) {
    $wsMain.Range($col + "2").Value2 = $edwinName
}

# --- Sheet 'eDWIN_LINK': rows 2-4 (description/name/purpose). Columns are
#     filled one language (column) at a time, top-to-bottom, to mirror the
#     original authoring order reflected in the shared-string table. ---
$wsLink = $wb.Worksheets.Item("eDWIN_LINK")
$wsLink.Range("C2").Value2 = $s24
$wsLink.Range("C3").Value2 = $s25
$wsLink.Range("C4").Value2 = $s26
$wsLink.Range("D2").Value2 = $s27
$wsLink.Range("D3").Value2 = $s28
$wsLink.Range("D4").Value2 = $s29
$wsLink.Range("E2").Value2 = $s30
$wsLink.Range("E3").Value2 = $s31
$wsLink.Range("E4").Value2 = $s32
$wsLink.Range("G2").Value2 = $s33
$wsLink.Range("G3").Value2 = $s34
$wsLink.Range("G4").Value2 = $s35
$wsLink.Range("H2").Value2 = $s36
$wsLink.Range("H3").Value2 = $s37
$wsLink.Range("H4").Value2 = $s38
$wsLink.Range("I2").Value2 = $s39
$wsLink.Range("I3").Value2 = $wsLink.Range("A3").Value2  # shared-string 19 (pl.gov.edwin.0_0_2.models.eDWIN_LINK.name) - reused verbatim per source data
$wsLink.Range("I4").Value2 = $s40
$wsLink.Range("J2").Value2 = $s41
$wsLink.Range("J3").Value2 = $s42
$wsLink.Range("J4").Value2 = $s43
$wsLink.Range("K2").Value2 = $s44
$wsLink.Range("K3").Value2 = $s45
$wsLink.Range("K4").Value2 = $s46
$wsLink.Range("L2").Value2 = $s47
$wsLink.Range("L3").Value2 = $s48
$wsLink.Range("L4").Value2 = $s49
$wsLink.Range("M2").Value2 = $s50
$wsLink.Range("M3").Value2 = $s51
$wsLink.Range("M4").Value2 = $s52

# --- Normalize row heights: multi-line descriptions can trigger Excel's
#     automatic row-height adjustment; AutoFit brings rows that don't
#     actually need extra height back to the sheet's default (no explicit
#     height / customHeight attribute), matching the original formatting. ---
$wsLink.Rows.Item(2).AutoFit() | Out-Null
$wsLink.Rows.Item(3).AutoFit() | Out-Null
$wsLink.Rows.Item(4).AutoFit() | Out-Null
$wsMain.Rows.Item(2).AutoFit() | Out-Null
